$wb = $excel.ActiveWorkbook

# --- "ir" sheet: Table7 row 2 (tradingPartner/createdBy) ---
$irWs = $wb.Worksheets.Item("ir")
$irWs.Range("A2").Value = "ARCONICTP"
$irWs.Range("B2").ClearContents()

# --- "tier" sheet: Table9 gains a second data row (append missing code before save rule) ---
$tierWs = $wb.Worksheets.Item("tier")
$tierTable = $tierWs.ListObjects.Item("Table9")
[void]$tierTable.ListRows.Add()

# Update the tradingPartner on the existing row to the new partner name.
$tierWs.Range("A2").Value = "ARCONICTP"

# Copy the rest of row 2's formatting/values down into row 3 (t1bt..t2cu stay identical),
# then fill in the two unique columns for the newly appended row.
$tierWs.Range("C2:L2").Copy($tierWs.Range("C3:L3"))
$tierWs.Range("A3").Value = "ARCONICTP"
$tierWs.Range("B3").Value = "ARCONICTP_HJBT 06.15.2020"

# Column B ("fileName") widens to fit the longer file name now stored in it.
$tierWs.Columns.Item(2).ColumnWidth = 34

# --- "simpleton" sheet: Table10 row 2 (fileName/tabName) gets renamed ---
$simpletonWs = $wb.Worksheets.Item("simpleton")
$simpletonWs.Range("A2").Value = "simpletonTest"
$simpletonWs.Range("B2").Value = "TEST"

# A stray column-D resize on the simpleton sheet (same width, but now flagged custom).
[void]$simpletonWs.Columns.Item(4)
$simpletonWs.Columns.Item(4).ColumnWidth = $simpletonWs.Columns.Item(4).ColumnWidth

# --- View/selection bookkeeping to match the saved state ---
$tierWs.Activate()
[void]$tierWs.Range("A2:XFD2").Select()

$missingWs = $wb.Worksheets.Item("missingCode")
$missingWs.Activate()
[void]$missingWs.Range("A1:G2").Select()

# "simpleton" becomes the active/selected tab on save.
$simpletonWs.Activate()
[void]$simpletonWs.Range("G21:G22").Select()
